$d = $word.ActiveDocument

# Locate the unique "Работник" label, then find the following 24-underscore
# signature-line placeholder that immediately follows it (skipping the
# spacing run in between).
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$found = $rng.Find.Execute("Работник", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate 'Работник' label"
}
$afterLabel = $rng.End

# Search for the 24-underscore placeholder starting right after the label.
$rng2 = $d.Range($afterLabel, $d.Content.End)
$found2 = $rng2.Find.Execute("________________________", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the underscore placeholder after 'Работник'"
}

$insertAt = $rng2.End

# Insert a single trailing space right after the placeholder underscores.
$ins = $d.Range($insertAt, $insertAt)
$ins.InsertAfter(" ")

# Force the newly-typed space onto its own run (matching the template's
# run-per-formatting-change convention) by nudging a property through a
# real change before restoring the original inherited formatting.
$newRun = $d.Range($insertAt, $insertAt + 1)
$newRun.Font.Size = 99
$newRun.Font.Name = "Times New Roman"
$newRun.Font.Size = 10
$newRun.LanguageID = "ru-RU"
